# Natmi following Dr Hou advice
# Rebuild LR-pairs rows for Cntf-Cntfr across all 4 sending clusters (ECs, FAPs, M2, sCs)
# and 2 target clusters (FAPs, sCs) -> 8 data rows total (was 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cntf"
$ws.Range("C2").Value = "Cntfr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 1.249087666666667
$ws.Range("H2").Value = 3.747263
$ws.Range("I2").Value = 0.3661911890049113
$ws.Range("J2").Value = 0.3661911890049113
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 5.666771333333333
$ws.Range("N2").Value = 17.000314
$ws.Range("O2").Value = 0.9648745059153377
$ws.Range("P2").Value = 0.9648745059153376
$ws.Range("Q2").Value = 7.078294182286888
$ws.Range("R2").Value = 63.704647640582
$ws.Range("S2").Value = 0.3533285425616638
$ws.Range("T2").Value = 0.3533285425616638

# Row 3: ECs -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cntf"
$ws.Range("C3").Value = "Cntfr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 1.249087666666667
$ws.Range("H3").Value = 3.747263
$ws.Range("I3").Value = 0.3661911890049113
$ws.Range("J3").Value = 0.3661911890049113
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 0.2062943333333333
$ws.Range("N3").Value = 0.618883
$ws.Range("O3").Value = 0.03512549408466232
$ws.Range("P3").Value = 0.03512549408466231
$ws.Range("Q3").Value = 0.2576797074698889
$ws.Range("R3").Value = 2.319117367229
$ws.Range("S3").Value = 0.01286264644324747
$ws.Range("T3").Value = 0.01286264644324747

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cntf"
$ws.Range("C4").Value = "Cntfr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5813396666666667
$ws.Range("H4").Value = 1.744019
$ws.Range("I4").Value = 0.1704295618581233
$ws.Range("J4").Value = 0.1704295618581232
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 5.666771333333333
$ws.Range("N4").Value = 17.000314
$ws.Range("O4").Value = 0.9648745059153377
$ws.Range("P4").Value = 0.9648745059153376
$ws.Range("Q4").Value = 3.294318957996222
$ws.Range("R4").Value = 29.648870621966
$ws.Range("S4").Value = 0.1644431392912242
$ws.Range("T4").Value = 0.1644431392912241

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cntf"
$ws.Range("C5").Value = "Cntfr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2.0
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5813396666666667
$ws.Range("H5").Value = 1.744019
$ws.Range("I5").Value = 0.1704295618581233
$ws.Range("J5").Value = 0.1704295618581232
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.2062943333333333
$ws.Range("N5").Value = 0.618883
$ws.Range("O5").Value = 0.03512549408466232
$ws.Range("P5").Value = 0.03512549408466231
$ws.Range("Q5").Value = 0.1199270789752222
$ws.Range("R5").Value = 1.079343710777
$ws.Range("S5").Value = 0.0059864225668991
$ws.Range("T5").Value = 0.005986422566899098

# Row 6: M2 -> FAPs
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cntf"
$ws.Range("C6").Value = "Cntfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 0.704738
$ws.Range("H6").Value = 2.114214
$ws.Range("I6").Value = 0.20660587166442
$ws.Range("J6").Value = 0.20660587166442
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 5.666771333333333
$ws.Range("N6").Value = 17.000314
$ws.Range("O6").Value = 0.9648745059153377
$ws.Range("P6").Value = 0.9648745059153376
$ws.Range("Q6").Value = 3.993589095910667
$ws.Range("R6").Value = 35.942301863196
$ws.Range("S6").Value = 0.199348738341415
$ws.Range("T6").Value = 0.1993487383414149

# Row 7: M2 -> sCs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cntf"
$ws.Range("C7").Value = "Cntfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 0.704738
$ws.Range("H7").Value = 2.114214
$ws.Range("I7").Value = 0.20660587166442
$ws.Range("J7").Value = 0.20660587166442
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.2062943333333333
$ws.Range("N7").Value = 0.618883
$ws.Range("O7").Value = 0.03512549408466232
$ws.Range("P7").Value = 0.03512549408466231
$ws.Range("Q7").Value = 0.1453834558846667
$ws.Range("R7").Value = 1.308451102962
$ws.Range("S7").Value = 0.007257133323005089
$ws.Range("T7").Value = 0.007257133323005086

# Row 8: sCs -> FAPs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cntf"
$ws.Range("C8").Value = "Cntfr"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2.0
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.8758606666666667
$ws.Range("H8").Value = 2.627582
$ws.Range("I8").Value = 0.2567733774725455
$ws.Range("J8").Value = 0.2567733774725454
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 5.666771333333333
$ws.Range("N8").Value = 17.000314
$ws.Range("O8").Value = 0.9648745059153377
$ws.Range("P8").Value = 0.9648745059153376
$ws.Range("Q8").Value = 4.963302117860889
$ws.Range("R8").Value = 44.669719060748
$ws.Range("S8").Value = 0.2477540857210348
$ws.Range("T8").Value = 0.2477540857210347

# Row 9: sCs -> sCs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cntf"
$ws.Range("C9").Value = "Cntfr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2.0
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.8758606666666667
$ws.Range("H9").Value = 2.627582
$ws.Range("I9").Value = 0.2567733774725455
$ws.Range("J9").Value = 0.2567733774725454
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.2062943333333333
$ws.Range("N9").Value = 0.618883
$ws.Range("O9").Value = 0.03512549408466232
$ws.Range("P9").Value = 0.03512549408466231
$ws.Range("Q9").Value = 0.1806850923228889
$ws.Range("R9").Value = 1.626165830906
$ws.Range("S9").Value = 0.009019291751510661
$ws.Range("T9").Value = 0.009019291751510656

